$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# The "Source:" block (rows 70-73) gets a blank separator row inserted above
# the "SBS Main Indicators..." line, and the Eurostat URL line loses its
# hyperlink (and drops down below the existing blank separator row).
# ---------------------------------------------------------------------------

# Drop the hyperlink while it is still anchored at its original cell (A72) -
# hyperlink anchors in this engine do not follow a later row-insert.
$ws.Range("A72").Hyperlinks.Delete()

# Insert a new blank row above the "SBS Main Indicators..." row (old row 71).
# This pushes: 71->72 (text), 72->73 (url, now unlinked), 73->74 (old blank).
$ws.Rows.Item(71).Insert()

# Swap the (now unlinked) url text down past the pre-existing blank row, so
# the blank separator sits between the title line and the url line.
$urlText = $ws.Range("A73").Value()
$ws.Range("A73").Value = ""
$ws.Range("A74").Value = $urlText

# Restyle the (now blank) A73 to match its italic "source" siblings.
$ws.Range("A73").Font().Italic = $true

# ---------------------------------------------------------------------------
# The "National Institute of Statistics" / "SBS Eurostat" citation rows
# (old rows 76-79, now shifted to 77-80 because of the insert above) have
# their long citation lines replaced with a plain repeat of the title above
# them.
# ---------------------------------------------------------------------------

$ws.Range("A78").Value = $ws.Range("A77").Value()
$ws.Range("A80").Value = $ws.Range("A79").Value()
